$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HD105")

# --- Add TOTAL row (row 17) ---
$ws.Range("A17").Value = "TOTAL"
$ws.Range("A17").NumberFormat = "@"
$ws.Range("B17").Formula = "=SUM(B2:B16)"
$ws.Range("C17").Formula = "=SUM(C2:C16)"
$ws.Range("D17").Formula = "=SUM(D2:D16)"

# --- New header columns E/F/G ---
$ws.Range("E1").Value = "DEM %"
$ws.Range("F1").Value = "REP %"
$ws.Range("G1").Value = "MARGIN"

# --- Percentage / margin formulas for rows 2-17 ---
$ws.Range("E2:E17").Formula = "=B2/D2"
$ws.Range("F2:F17").Formula = "=C2/D2"
$ws.Range("G2:G17").Formula = "=(B2-C2)/D2"

$ws.Range("E2:G17").NumberFormat = "0.0%"

# --- Switch the active tab from "2018 Results" to "HD105" ---
$ws.Activate() | Out-Null
$ws.Range("H17").Select() | Out-Null
